$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "Summary": update aggregate metrics
# ---------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1200.24   # Current Capital
$summary.Range("B4").Value = 0.23      # Total P&L $
$summary.Range("B6").Value = 58        # Total Trades
$summary.Range("B7").Value = 25        # Winning Trades
$summary.Range("B9").Value = 43.1      # Win Rate %

# ---------------------------------------------------------------
# Sheet "Strategy Status": update the MarketMaking strategy row (row 4)
# ---------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 100.24     # Capital
$status.Range("D4").Value = 58         # Trades
$status.Range("E4").Value = 0.23       # P&L $
$status.Range("F4").Value = 0.24       # P&L %
$status.Range("G4").Value = 43.1       # Win Rate %

# ---------------------------------------------------------------
# Append the new closed trade (Trade #58) to both the "All Trades"
# and "MarketMaking" logs as row 59
# ---------------------------------------------------------------
$newRow = @(58, "2026-02-17", "12:52:27", "MarketMaking", "UP", 0.39, 0.4, "CLOSED", 2.5641, 0.01, 100.24, 0, 0, 0.6, "Normal spread capture: 19600 bps", "early_exit", 0.13)

foreach ($sheetName in @("All Trades", "MarketMaking")) {
    $ws = $wb.Worksheets.Item($sheetName)
    for ($col = 1; $col -le $newRow.Length; $col++) {
        $cell = $ws.Cells.Item(59, $col)
        $value = $newRow[$col - 1]
        if ($col -eq 2) {
            # Column B holds the date as plain text (e.g. "2026-02-17").
            # A bare assignment gets auto-recognised as a real date by
            # Excel's smart typing, so force text entry with a leading
            # apostrophe and then strip the resulting quote-prefix style
            # so the cell ends up as a plain, unstyled text cell.
            $cell.Value = "'" + $value
            $cell.ClearFormats()
        } else {
            $cell.Value = $value
        }
    }
}
